# "Foi implementado o sistema de edição selecionável em todas as listas"
# Adds the historical tracking rows (2-5) to the DERMO list sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new rows hold values that look numeric (dates, money, percentages)
# but must be stored as literal text, exactly as typed — format the
# target range as Text first so Excel doesn't auto-convert them to
# numbers/dates.
$ws.Range("A2:H5").NumberFormat = "@"

$data = @(
    @("31/03/2001", "3000.00", "3000.00", "3000.00", "3000.00", "20", "0.00", "100.00"),
    @("08/07/2023", "4000.00", "7000.00", "4000.00", "7000.00", "40", "0.00", "100.00"),
    @("08/07/2023", "1000.00", "8000.00", "1000.00", "8000.00", "50", "0.00", "100.00"),
    @("08/07/2023", "4000.00", "12000.00", "1000.00", "9000.00", "60", "3000.00", "75.00")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $row[$j]
    }
}
